# Daily attendance processing - 2026-01-15 23:37:48
# Swap the order of "System" and the email address in the "Recorded By"
# column (G) wherever both values are present together as
# "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$used = $ws.UsedRange
$startRow = $used.Row
$rowCount = $used.Rows.Count
$lastRow = $startRow + $rowCount - 1

$changed = 0
for ($r = $startRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()
    if ($val -eq $oldValue) {
        $cell.Value = $newValue
        $changed = $changed + 1
    }
}

Write-Host "Updated $changed cell(s) in column G from '$oldValue' to '$newValue'."
